# FlixelRL-190 Floor1-3: enemy strength tuning ("player" sheet).
# - Widen column B (exp) to fit the much larger numbers now stored there.
# - Lower the lv1 hp/str (C3/D3).
# - Rescale the exp column (B4:B102) to the new, much steeper curve.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("exp") width: 5.13281 -> 7.21875
$ws.Columns.Item(2).ColumnWidth = 7.21875

# lv1 (row 3) hp/str tweak
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 3

# New exp requirements for lv2..lv100 (rows 4..102), column B
$newExp = @(
    20,50,100,250,500,800,1200,1800,3000,5000,
    10000,30000,50000,80000,120000,170000,250000,500000,1000000,2000000,
    3000000,4000000,5000000,6000000,7000000,8000000,9000000,10000000,11000000,12000000,
    13000000,14000000,15000000,16000000,17000000,18000000,19000000,20000000,21000000,22000000,
    23000000,24000000,25000000,26000000,27000000,28000000,29000000,30000000,31000000,32000000,
    33000000,34000000,35000000,36000000,37000000,38000000,39000000,40000000,41000000,42000000,
    43000000,44000000,45000000,46000000,47000000,48000000,49000000,50000000,51000000,52000000,
    53000000,54000000,55000000,56000000,57000000,58000000,59000000,60000000,61000000,62000000,
    63000000,64000000,65000000,66000000,67000000,68000000,69000000,70000000,71000000,72000000,
    73000000,74000000,75000000,76000000,77000000,78000000,79000000,80000000,81000000
)

$startRow = 4
for ($i = 0; $i -lt $newExp.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $newExp[$i]
}
